$wb = $excel.ActiveWorkbook

# 1. Rename the translation sheet: "TRN-1079-0833-9890 (FA)" -> "FA (TRN-1079-0833-9890)"
$ws = $wb.Worksheets.Item("TRN-1079-0833-9890 (FA)")
$ws.Name = "FA (TRN-1079-0833-9890)"

# 2. Fix the header row of the translation attributes sheet.
#    Old headers: Key | Original Value | Action | Value | Comment
#    New headers: key | <Locale> (TRN-...) | action | value | comment | editor
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "Persian (TRN-1079-0833-9890)"
$ws.Range("C1").Value = "action"
$ws.Range("D1").Value = "value"
$ws.Range("E1").Value = "comment"
$ws.Range("F1").Value = "editor"

# 3. Populate the new "editor" column for every data row with the default "-" value.
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("F$r").Value = "-"
}

# 4. Move the active selection like in the saved workbook.
$ws.Activate()
$ws.Range("E18").Select()
